$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 2..15: column B (id) and column C (speaker_variant).
# Column D (is_prefered) is cleared for every data row.
$rows = @(
    @{ B = "#laodice";  C = "Laodice" },
    @{ B = "#axiane";   C = "Axiane" },
    @{ B = "#theóxena"; C = "Theóxena" },
    @{ B = "#cleomenes";C = "Cleomenes" },
    @{ B = "#atis";     C = "Atis" },
    @{ B = "#theoxena"; C = "Theoxena" },
    @{ B = "#theoxfna"; C = "Theoxfna" },
    @{ B = "#evander";  C = "Evander" },
    @{ B = "#cleovenes";C = "Cleovenes" },
    @{ B = "#ariarates";C = "Ariarates" },
    @{ B = "#attalus";  C = "Attalus" },
    @{ B = "#gentius";  C = "Gentius" },
    @{ B = "#dorine";   C = "Dorine" },
    @{ B = "#didos";    C = "Didos" }
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $rows[$i].B
    $ws.Cells.Item($r, 3).Value = $rows[$i].C
    $ws.Cells.Item($r, 4).ClearContents()
}
